$d = $word.ActiveDocument

$replacements = @(
    @{old="928÷6="; new="429÷2="},
    @{old="458÷5="; new="942÷2="},
    @{old="560÷9="; new="646÷8="},
    @{old="307÷4="; new="740÷5="},
    @{old="360÷6="; new="351÷5="},
    @{old="166÷4="; new="616÷3="},
    @{old="730÷7="; new="597÷3="},
    @{old="108÷2="; new="313÷2="},
    @{old="919÷8="; new="503÷3="},
    @{old="757÷2="; new="143÷2="},
    @{old="163÷3="; new="725÷2="},
    @{old="683÷2="; new="648÷5="},
    @{old="873÷8="; new="932÷4="},
    @{old="816÷6="; new="442÷9="},
    @{old="324÷9="; new="593÷8="},
    @{old="486÷3="; new="311÷5="},
    @{old="208÷7="; new="975÷3="},
    @{old="122÷7="; new="555÷2="},
    @{old="116÷2="; new="377÷6="},
    @{old="345÷3="; new="411÷5="},
    @{old="349÷7="; new="337÷9="},
    @{old="397÷8="; new="441÷3="},
    @{old="713÷6="; new="758÷5="},
    @{old="770÷4="; new="214÷9="},
    @{old="252÷7="; new="299÷6="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
